# Reorder data rows 2-13 (columns A:AY) according to the mapping below.
# Mapping key = destination row (after), value = source row (before),
# i.e. the data currently sitting in the "source" row should end up in
# the "destination" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRange = $ws.Range("A2:AY13")
$values = $srcRange.Value2

# Columns Y and AA hold plain-text dates (e.g. "2023-09-11"). Force the
# cells to Text format before writing back so Excel does not silently
# convert the strings into date serial numbers.
$ws.Range("Y2:Y13").NumberFormat = "@"
$ws.Range("AA2:AA13").NumberFormat = "@"

# Mapping: destination row offset (1-based within the 2..13 block) -> source row offset
$mapping = @{
    1  = 4   # row 2  <- row 5
    2  = 11  # row 3  <- row 12
    3  = 6   # row 4  <- row 7
    4  = 2   # row 5  <- row 3
    5  = 9   # row 6  <- row 10
    6  = 8   # row 7  <- row 9
    7  = 7   # row 8  <- row 8
    8  = 12  # row 9  <- row 13
    9  = 3   # row 10 <- row 4
    10 = 10  # row 11 <- row 11
    11 = 5   # row 12 <- row 6
    12 = 1   # row 13 <- row 2
}

$numCols = 51  # columns A..AY
$newValues = New-Object 'object[,]' 12, $numCols

for ($destOffset = 1; $destOffset -le 12; $destOffset++) {
    $srcOffset = $mapping[$destOffset]
    for ($col = 1; $col -le $numCols; $col++) {
        $newValues[$destOffset - 1, $col - 1] = $values[$srcOffset, $col]
    }
}

$srcRange.Value2 = $newValues
